$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new "site group" header rows (site name only, in column A) ---
# Final row 5: "Meduxnekeag Lake"
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Meduxnekeag Lake"
$ws.Range("B5:F5").Clear()

# Final row 8: "Fraser"
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Fraser"
$ws.Range("B8:F8").Clear()

# Final row 18: "Pontoosuc Lake"
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "Pontoosuc Lake"
$ws.Range("B18:F18").Clear()

# --- Add the grand-total row at the bottom (row 37) ---
$ws.Range("D37").Value = "Total:"
$ws.Range("E37").Formula = "=SUM(E2:E35)"

Write-Output "edit complete"
